# ------------------------------------------------------------------
# Adds a new "Player Info" sheet (as the first sheet) with player
# biographical data, and changes the MATCH_CARD_LINK column (full URL)
# into a MATCH_CODE column (just the numeric match code) on both the
# "ODI Batting" and "ODI Bowling" sheets.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "Player Info" worksheet. Worksheets.Add() inserts
#    the new sheet before the currently active sheet, which places it
#    first (matching the target sheet order: Player Info, ODI Batting,
#    ODI Bowling).
#    NOTE: sheet references returned by Worksheets.Item(...) appear to
#    be positional, so fetch the "ODI Batting"/"ODI Bowling" sheet
#    objects *after* the new sheet has been inserted, not before.
# ------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$playerInfo.Range("A2").Value = "'3501"
$playerInfo.Range("B2").Value = "Mahawaduge Dilruwan Kamalaneth Perera"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# Copy the bold/centered/bordered header style used by the other
# sheets' header rows onto the new header row.
$battingSheet.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# ------------------------------------------------------------------
# 2. "ODI Batting" sheet: rename MATCH_CARD_LINK -> MATCH_CODE (D1),
#    and replace the full scorecard URLs in D2:D14 with just the
#    trailing numeric match code.
# ------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingCodes = @{
    2  = "2741"
    3  = "2795"
    4  = "2798"
    5  = "2799"
    6  = "3718"
    7  = "3923"
    8  = "3924"
    9  = "3927"
    10 = "3929"
    11 = "3931"
    12 = "4012"
    13 = "4014"
    14 = "4193"
}

foreach ($row in $battingCodes.Keys) {
    $battingSheet.Range("D$row").Value = "'" + $battingCodes[$row]
}

# ------------------------------------------------------------------
# 3. "ODI Bowling" sheet: rename MATCH_CARD_LINK -> MATCH_CODE (B1),
#    and replace the full scorecard URLs in B2:B11 with just the
#    trailing numeric match code.
# ------------------------------------------------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @{
    2  = "2798"
    3  = "2799"
    4  = "3718"
    5  = "3923"
    6  = "3924"
    7  = "3927"
    8  = "3929"
    9  = "3931"
    10 = "4014"
    11 = "4193"
}

foreach ($row in $bowlingCodes.Keys) {
    $bowlingSheet.Range("B$row").Value = "'" + $bowlingCodes[$row]
}
